$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New job posting: Job_Id=15 (row 16), Jd_Title = "Cyber Security Engineer",
# Job_Description = "Demo"
$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Cyber Security Engineer"
$ws.Range("C16").Value = "Demo"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
